$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row/column -> new text for the cells whose stats changed.
$updates = @(
    @{ Row = 2; Col = 2; Old = "305 (100.0)"; New = "306 (100.0)" },
    @{ Row = 2; Col = 3; Old = "279 (100.0)"; New = "280 (100.0)" },
    @{ Row = 3; Col = 2; Old = "279 (91.8)";  New = "280 (91.8)" },
    @{ Row = 3; Col = 3; Old = "279 (100.0)"; New = "280 (100.0)" },
    @{ Row = 4; Col = 2; Old = "204 (69.4)";  New = "205 (69.5)" },
    @{ Row = 4; Col = 3; Old = "204 (75.3)";  New = "205 (75.4)" },
    @{ Row = 5; Col = 2; Old = "156 (55.9)";  New = "157 (55.9)" },
    @{ Row = 5; Col = 3; Old = "156 (60.2)";  New = "157 (60.2)" },
    @{ Row = 6; Col = 2; Old = "142 (53.8)";  New = "143 (54.0)" },
    @{ Row = 6; Col = 3; Old = "142 (57.7)";  New = "143 (57.9)" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    $current = $rng.Text.TrimEnd([char]13, [char]7)
    if ($current -ne $u.Old) {
        Write-Output ("MISMATCH at row " + $u.Row + " col " + $u.Col + ": expected '" + $u.Old + "' found '" + $current + "'")
    }
    $rng.Text = $u.New
}

Write-Output "Done"
